$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,(37.0, 'Std_RF_DummyEncoded_MinMaxScaling', 0.91, 0.6, 0.27, 0.37, 9.992, 1.0, 3.3966, 0.91)
  ,(38.0, 'Std_RF_DummyEncoded_MinMaxScaling_SMOTE', 0.9, 0.43, 0.42, 0.42, 9.992, 1.0, 4.2363, 0.8773)
  ,(39.0, 'Std_RF_DummyEncoded_SMOTE', 0.91, 0.57, 0.32, 0.41, 9.992, 1.0, 3.4348, 0.9005)
  ,(40.0, 'Std_RF_DummyEncoded', 0.91, 0.6, 0.26, 0.36, 9.992, 1.0, 3.4348, 0.9005)
  ,(41.0, 'Std_RF_DummyEncoded_Binning', 0.81, 0.5, 0.28, 0.36, 0.5444, 0.9842, 3.7401, 0.8917)
  ,(42.0, 'Std_RF_DummyEncoded_Binning_SMOTE', 0.82, 0.44, 0.36, 0.4, 0.3137, 0.9909, 4.0836, 0.8818)
  ,(43.0, 'Std_RF_LabeEncoded_MinMaxScaling', 0.91, 0.61, 0.28, 0.38, 9.992, 1.0, 3.3585, 0.9028)
  ,(44.0, 'Std_RF_LabelEncoded_MinMaxScaling_SMOTE', 0.9, 0.4, 0.51, 0.45, 9.992, 1.0, 4.6943, 0.8641)
  ,(45.0, 'Std_RF_LabelEncoded_SMOTE', 0.89, 0.45, 0.54, 0.49, 9.992, 1.0, 4.2363, 0.8773)
  ,(46.0, 'Std_RF_LabelEncoded', 0.91, 0.6, 0.27, 0.37, 9.992, 1.0, 3.3966, 0.9016)
  ,(47.0, 'Std_RF_LabelEncoded_Binning', 0.8, 0.47, 0.28, 0.35, 0.5444, 0.9842, 3.8928, 0.8873)
  ,(48.0, 'Std_RF_LabelEncoded_Binning_SMOTE', 0.81, 0.37, 0.4, 0.38, 0.4814, 0.9861, 4.7706, 0.8619)
  ,(49.0, 'RS_RF_DummyEncoded_MinMaxScaling', 0.92, 0.53, 0.17, 0.26, 1.3086, 0.9621, 3.6638, 0.8939)
  ,(50.0, 'RS_RF_DummyEncoded_MinMaxScaling_SMOTE', 0.91, 0.48, 0.42, 0.45, 0.0054, 0.9998, 3.8546, 0.8884)
  ,(51.0, 'RS_RF_DummyEncoded_SMOTE', 0.92, 0.53, 0.32, 0.4, 0.2921, 0.9915, 3.5875, 0.8961)
  ,(52.0, 'RS_RF_DummyEncoded', 0.92, 0.53, 0.17, 0.26, 1.3086, 0.9621, 3.6638, 0.8939)
  ,(53.0, 'RS_RF_DummyEncoded_Binning', 0.83, 0.55, 0.21, 0.31, 3.0088, 0.9129, 3.5875, 0.8961)
  ,(54.0, 'RS_RF_DummyEncoded_Binning_SMOTE', 0.82, 0.48, 0.37, 0.42, 0.7518, 0.9782, 3.8546, 0.8884)
  ,(55.0, 'RS_RF_LabeEncoded_MinMaxScaling', 0.92, 0.59, 0.24, 0.35, 1.8148, 0.9474, 3.473, 0.8994)
  ,(56.0, 'RS_RF_LabelEncoded_MinMaxScaling_SMOTE', 0.9, 0.44, 0.52, 0.47, 9.992, 1.0, 4.3126, 0.8751)
  ,(57.0, 'RS_RF_LabelEncoded_SMOTE', 0.9, 0.43, 0.51, 0.47, 9.992, 1.0, 4.389, 0.8729)
  ,(58.0, 'RS_RF_LabelEncoded', 0.92, 0.59, 0.24, 0.35, 1.8148, 0.9474, 3.473, 0.8994)
  ,(59.0, 'RS_RF_LabelEncoded_Binning', 0.81, 0.53, 0.23, 0.33, 2.407, 0.9303, 3.6256, 0.895)
  ,(60.0, 'RS_RF_LabelEncoded_Binning_SMOTE', 0.81, 0.43, 0.41, 0.42, 0.7356, 0.9787, 4.2363, 0.8773)
)

$startRow = 40
for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $startRow + $i
  $row = $data[$i]
  for ($c = 0; $c -lt $row.Count; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
}

$ws.Range("J1").Select() | Out-Null
